$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FE)
$ws.Range("C2").Value = 0.83
$ws.Range("D2").Value = 1.85
$ws.Range("E2").Value = 0.21

# Row 3 (FE+Disg)
$ws.Range("C3").Value = 0.83
$ws.Range("D3").Value = 1.85
$ws.Range("E3").Value = 0.2
$ws.Range("F3").Value = 0.98
$ws.Range("G3").Value = 0.1

# Row 4 (FE+Disg+Var)
$ws.Range("B4").Value = 0.35
$ws.Range("C4").Value = 0.83
$ws.Range("D4").Value = 1.85
$ws.Range("E4").Value = 0.39
$ws.Range("G4").Value = 0.77
